$d = $word.ActiveDocument

# Locate the paragraph that holds item "29. Aging and Mental Health"
$find = $d.Content
$find.Find.ClearFormatting()
$null = $find.Find.Execute("Aging and Mental Health", $false, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
$para29 = $find.Paragraphs(1)

# Insert a brand-new paragraph right after it (before the blank paragraph that follows)
$endRange = $para29.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$newPara = $para29.Next()

# Add "30. " as its own run
$newPara.Range.InsertAfter("30. ")

# Split the new paragraph into two (so the "30. " text and the title end up in
# separate runs, matching the source document's pattern for every other entry),
# then immediately merge the two paragraphs back together by deleting the
# paragraph mark that separated them. The two runs stay distinct even after
# the merge, and neither run picks up any stray character formatting.
$splitPoint = $newPara.Range.End - 1
$splitRange = $d.Range($splitPoint, $splitPoint)
$splitRange.InsertParagraphAfter()

$titlePara = $newPara.Next()
$titlePara.Range.InsertAfter("International Journal of Geriatric Psychiatry")

$mergeMark = $newPara.Range.End - 1
$d.Range($mergeMark, $mergeMark + 1).Delete()

# Match the first-line indent used by every other entry in the list (720 twips = 36pt)
$newPara.Range.ParagraphFormat.FirstLineIndent = 36
